$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.628.92'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '2.268.01'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '119.26'
$ws.Range('E5').Value = '  +8.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '265.49'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.643'
$ws.Range('E7').Value = '  +4.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.01'
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  +3.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.57'
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0942'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.55'
$ws.Range('E12').Value = '  +9.60%  '
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.57'
$ws.Range('E14').Value = '  +2.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.899'
$ws.Range('E15').Value = '  +6.44%  '
$ws.Range('D16').Value = '2.612.87'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').Value = '2.268.75'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '43.562.60'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').Value = '  +2.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.88'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.15'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.41'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.07'
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.55'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.89'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.02'
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.84'
$ws.Range('E27').Value = '  +6.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.13'
$ws.Range('E28').Value = '  +5.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.37'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.24'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.34'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.53'
$ws.Range('E32').Value = '  +2.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0917'
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.73'
$ws.Range('E34').Value = '  +0.30%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.34'
$ws.Range('E35').Value = '  +14.06%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.130'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0389'
$ws.Range('E37').Value = '  +11.67%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.112'
$ws.Range('E38').Value = '  +8.43%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.63'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.53'
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.69'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.59'
$ws.Range('E42').Value = '  -4.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.237'
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.39'
$ws.Range('E45').Value = '  +3.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.78'
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '80.39'
$ws.Range('E47').Value = '  +54.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.27'
$ws.Range('E48').Value = '  +4.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.88'
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.476'
$ws.Range('E51').Value = '  +9.26%  '
